# [Kadastro App] Yeni kayit eklendi: 3022
# Adds a new record row (row 81) to both the "Kayitlar" (all records) sheet
# and the "Erdemli" (district-filtered) sheet, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$newRow = @("3022", "2025-09-11", "Erdemli", "1", "3B", "EMİNE ALANLI KIRCILI (K.Mühendisi), SERDAR ARSLAN (Tekniker)")

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Make sure the new cells are stored as text (matching the existing
    # numberStoredAsText rows) rather than being auto-coerced to numbers
    # or dates by Excel's normal Value-assignment heuristics.
    $targetRange = $ws.Range("A81:F81")
    $targetRange.NumberFormat = "@"

    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item(81, $col).Value = $newRow[$col - 1]
    }
}
